$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 386 (pushes existing rows 386-418 down to 387-419)
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row with its own data
$ws.Cells.Item(386, 1).Value = 4
$ws.Cells.Item(386, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(386, 3).Value = "Los Lagos"
$ws.Cells.Item(386, 4).Value = 45223
$ws.Cells.Item(386, 5).Value = 10
$ws.Cells.Item(386, 6).Value = 100112039
$ws.Cells.Item(386, 7).Value = "Ciboulette"
$ws.Cells.Item(386, 8).Value = "Sin especificar"
$ws.Cells.Item(386, 9).Value = "Primera"
$ws.Cells.Item(386, 10).Value = 240
$ws.Cells.Item(386, 11).Value = 3500
$ws.Cells.Item(386, 12).Value = 3500
$ws.Cells.Item(386, 13).Value = 3500
$ws.Cells.Item(386, 14).Value = "`$/docena de atados"
$ws.Cells.Item(386, 15).Value = "Región Metropolitana"
$ws.Cells.Item(386, 16).Value = 1167
$ws.Cells.Item(386, 17).Value = 3
$ws.Cells.Item(386, 18).Value = "Hortaliza"
